{"js": "const replacements = [\n  { oldText: \"698\u00f77=99, 5\",  newText: \"887\u00f79=98, 5\" },\n  { oldText: \"163\u00f74=40, 3\",  newText: \"755\u00f74=188, 3\" },\n  { oldText: \"480\u00f72=240, 0\", newText: \"656\u00f72=328, 0\" },\n  { oldText: \"144\u00f79=16, 0\",  newText: \"193\u00f76=32, 1\" },\n  { oldText: \"804\u00f74=201, 0\", newText: \"363\u00f74=90, 3\" },\n  { oldText: \"112\u00f78=14, 0\",  newText: \"215\u00f75=43, 0\" },\n  { oldText: \"765\u00f75=153, 0\", newText: \"352\u00f79=39, 1\" },\n  { oldText: \"872\u00f75=174, 2\", newText: \"165\u00f77=23, 4\" },\n  { oldText: \"565\u00f75=113, 0\", newText: \"252\u00f78=31, 4\" },\n  { oldText: \"983\u00f72=491, 1\", newText: \"615\u00f73=205, 0\" },\n  { oldText: \"656\u00f73=218, 2\", newText: \"251\u00f74=62, 3\" },\n  { oldText: \"416\u00f73=138, 2\", newText: \"319\u00f72=159, 1\" },\n  { oldText: \"308\u00f78=38, 4\",  newText: \"700\u00f72=350, 0\" },\n  { oldText: \"674\u00f74=168, 2\", newText: \"428\u00f74=107, 0\" },\n  { oldText: \"885\u00f79=98, 3\",  newText: \"580\u00f77=82, 6\" },\n  { oldText: \"615\u00f74=153, 3\", newText: \"777\u00f78=97, 1\" },\n  { oldText: \"924\u00f72=462, 0\", newText: \"803\u00f77=114, 5\" },\n  { oldText: \"427\u00f77=61, 0\",  newText: \"603\u00f73=201, 0\" },\n  { oldText: \"923\u00f76=153, 5\", newText: \"220\u00f76=36, 4\" },\n  { oldText: \"536\u00f75=107, 1\", newText: \"711\u00f76=118, 3\" },\n  { oldText: \"645\u00f79=71, 6\",  newText: \"787\u00f79=87, 4\" },\n  { oldText: \"955\u00f72=477, 1\", newText: \"261\u00f73=87, 0\" },\n  { oldText: \"511\u00f75=102, 1\", newText: \"678\u00f78=84, 6\" },\n  { oldText: \"154\u00f77=22, 0\",  newText: \"535\u00f75=107, 0\" },\n  { oldText: \"202\u00f74=50, 2\",  newText: \"396\u00f77=56, 4\" }\n];\n\nconst body = context.document.body;\n\nfor (const { oldText, newText } of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"698\u00f77=99, 5\";  New = \"887\u00f79=98, 5\" },\n    @{ Old = \"163\u00f74=40, 3\";  New = \"755\u00f74=188, 3\" },\n    @{ Old = \"480\u00f72=240, 0\"; New = \"656\u00f72=328, 0\" },\n    @{ Old = \"144\u00f79=16, 0\";  New = \"193\u00f76=32, 1\" },\n    @{ Old = \"804\u00f74=201, 0\"; New = \"363\u00f74=90, 3\" },\n    @{ Old = \"112\u00f78=14, 0\";  New = \"215\u00f75=43, 0\" },\n    @{ Old = \"765\u00f75=153, 0\"; New = \"352\u00f79=39, 1\" },\n    @{ Old = \"872\u00f75=174, 2\"; New = \"165\u00f77=23, 4\" },\n    @{ Old = \"565\u00f75=113, 0\"; New = \"252\u00f78=31, 4\" },\n    @{ Old = \"983\u00f72=491, 1\"; New = \"615\u00f73=205, 0\" },\n    @{ Old = \"656\u00f73=218, 2\"; New = \"251\u00f74=62, 3\" },\n    @{ Old = \"416\u00f73=138, 2\"; New = \"319\u00f72=159, 1\" },\n    @{ Old = \"308\u00f78=38, 4\";  New = \"700\u00f72=350, 0\" },\n    @{ Old = \"674\u00f74=168, 2\"; New = \"428\u00f74=107, 0\" },\n    @{ Old = \"885\u00f79=98, 3\";  New = \"580\u00f77=82, 6\" },\n    @{ Old = \"615\u00f74=153, 3\"; New = \"777\u00f78=97, 1\" },\n    @{ Old = \"924\u00f72=462, 0\"; New = \"803\u00f77=114, 5\" },\n    @{ Old = \"427\u00f77=61, 0\";  New = \"603\u00f73=201, 0\" },\n    @{ Old = \"923\u00f76=153, 5\"; New = \"220\u00f76=36, 4\" },\n    @{ Old = \"536\u00f75=107, 1\"; New = \"711\u00f76=118, 3\" },\n    @{ Old = \"645\u00f79=71, 6\";  New = \"787\u00f79=87, 4\" },\n    @{ Old = \"955\u00f72=477, 1\"; New = \"261\u00f73=87, 0\" },\n    @{ Old = \"511\u00f75=102, 1\"; New = \"678\u00f78=84, 6\" },\n    @{ Old = \"154\u00f77=22, 0\";  New = \"535\u00f75=107, 0\" },\n    @{ Old = \"202\u00f74=50, 2\";  New = \"396\u00f77=56, 4\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($r.Old, $false, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
